$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @("华胜天成", "华胜天成", "华胜天成")
    3  = @("航天发展", "航天发展", "豫能控股")
    4  = @("润泽科技", "中天科技", "协鑫集成")
    5  = @("利欧股份", "利欧股份", "沪电股份")
    6  = @("云南锗业", "润泽科技", "利欧股份")
    7  = @("协鑫集成", "包钢股份", "博纳影业")
    8  = @("胜宏科技", "云南锗业", "润泽科技")
    9  = @("中天科技", "岩山科技", "航天发展")
    10 = @("汉缆股份", "协鑫集成", "汉缆股份")
    11 = @("金风科技", "汉缆股份", "大位科技")
    12 = @("大位科技", "沪电股份", "明阳电路")
    13 = @("包钢股份", "亨通光电", "金正大")
    14 = @("亨通光电", "豫能控股", "金风科技")
    15 = @("岩山科技", "金正大", "中天科技")
    16 = @("豫能控股", "金风科技", "洲际油气")
    17 = @("泰嘉股份", "聚飞光电", "章源钨业")
    18 = @("川润股份", "章源钨业", "杭电股份")
    19 = @("掌阅科技", "东方电气", "云南锗业")
    20 = @("金正大", "航天动力", "大族激光")
    21 = @("华丰科技", "明阳电路", "华银电力")
}

foreach ($rowNum in $data.Keys) {
    $values = $data[$rowNum]
    $ws.Cells.Item($rowNum, 1).Value = $values[0]
    $ws.Cells.Item($rowNum, 2).Value = $values[1]
    $ws.Cells.Item($rowNum, 3).Value = $values[2]
}
